$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph so we can
# insert the new bullet-list paragraph with the two professors right
# after it (and before "Programa resumido").
$targetPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "Docente(s) Responsável(eis)*") {
        $targetPara = $para
        break
    }
}

if ($targetPara -eq $null) {
    throw "Could not find the 'Docente(s) Responsável(eis)' paragraph"
}

# Create a brand-new empty paragraph right after it.
$targetPara.Range.InsertParagraphAfter()

# Grab that freshly minted (currently empty) paragraph and inject the
# exact OOXML for the two professors, separated by a manual line break,
# styled as a bulleted list item - matching the target markup exactly:
#   <w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr>
#     <w:r><w:t>471420 - Carlos Antonio Reis Pereira Baptista</w:t><w:br/></w:r>
#     <w:r><w:t>3586455 - Cassius Olivio Figueiredo Terra Ruchert</w:t></w:r>
#   </w:p>
$newPara = $targetPara.Next()
$newRange = $newPara.Range

$ooxml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>471420 - Carlos Antonio Reis Pereira Baptista</w:t><w:br/></w:r><w:r><w:t>3586455 - Cassius Olivio Figueiredo Terra Ruchert</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$newRange.InsertXML($ooxml)

Write-Host "Inserted Docente(s) paragraph with 2 professors."
